$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$excel.Goto($ws.Range("C5"), $true)
